$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date header
$ws.Range("B1").Value = "'07/04/2023"

# Updated hour (A) / value (B) pairs, shifted by one hour (index starts at 9)
$data = @(
    @(9, 94),
    @(10, 200),
    @(11, 236),
    @(12, 214),
    @(13, 168),
    @(14, 133),
    @(15, 129),
    @(16, 128),
    @(17, 132),
    @(18, 135),
    @(19, 132),
    @(20, 104),
    @(21, 73),
    @(22, 43),
    @(23, 22),
    @(24, 10)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
